# Append additional training data rows (21-50) to the TrainData worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(21, 21, 0.05214, "semi-critical"),
    @(22, 22, 0.02215, "non-critical"),
    @(23, 23, 0.15564, "critical"),
    @(24, 24, 0.18546, "critical"),
    @(25, 25, 0.07456, "semi-critical"),
    @(26, 26, 0.08431, "semi-critical"),
    @(27, 27, 0.21543, "critical"),
    @(28, 28, 0.00512, "non-critical"),
    @(29, 29, 0.04152, "non-critical"),
    @(30, 30, 0.12453, "critical"),
    @(31, 31, 0.06445, "semi-critical"),
    @(32, 32, 0.07754, "semi-critical"),
    @(33, 33, 0.09454, "semi-critical"),
    @(34, 34, 0.04225, "non-critical"),
    @(35, 35, 0.03445, "non-critical"),
    @(36, 36, 0.05231, "non-critical"),
    @(37, 37, 0.03125, "non-critical"),
    @(38, 38, 0.07335, "semi-critical"),
    @(39, 39, 0.24156, "critical"),
    @(40, 40, 0.09556, "semi-critical"),
    @(41, 41, 0.12412, "critical"),
    @(42, 42, 0.01521, "non-critical"),
    @(43, 43, 0.04114, "non-critical"),
    @(44, 44, 0.05124, "non-critical"),
    @(45, 45, 0.09024, "semi-critical"),
    @(46, 46, 0.34152, "critical"),
    @(47, 47, 0.07541, "semi-critical"),
    @(48, 48, 0.02415, "non-critical"),
    @(49, 49, 0.04215, "non-critical"),
    @(50, 50, 0.06124, "semi-critical")
)

$row = 22
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# Scroll the viewport so row 29 is at the top (matches the author's
# on-screen position after typing the new rows), then leave the
# selection where editing stopped, one row below the last data row.
$excel.Goto($ws.Range("A29"), $true)
$ws.Range("D52").Select()
